$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 901
$ws.Cells.Item(5, 6).Value = 1073
$ws.Cells.Item(6, 6).Value = 1548
$ws.Cells.Item(7, 6).Value = 330
$ws.Cells.Item(8, 6).Value = 671
$ws.Cells.Item(9, 6).Value = 12146
$ws.Cells.Item(10, 6).Value = 125
$ws.Cells.Item(11, 6).Value = 2164
$ws.Cells.Item(13, 6).Value = 254
$ws.Cells.Item(15, 6).Value = 1224
$ws.Cells.Item(16, 6).Value = 204
$ws.Cells.Item(17, 6).Value = 271
$ws.Cells.Item(18, 6).Value = 773
$ws.Cells.Item(19, 6).Value = 674
$ws.Cells.Item(20, 6).Value = 296
$ws.Cells.Item(22, 6).Value = 754
$ws.Cells.Item(23, 6).Value = 4065
$ws.Cells.Item(24, 6).Value = 1114
$ws.Cells.Item(25, 6).Value = 860
$ws.Cells.Item(29, 6).Value = 1045
$ws.Cells.Item(30, 6).Value = 48
$ws.Cells.Item(31, 6).Value = 98
$ws.Cells.Item(33, 6).Value = 28
$ws.Cells.Item(35, 6).Value = 27
$ws.Cells.Item(36, 6).Value = 12
$ws.Cells.Item(37, 6).Value = 4415
$ws.Cells.Item(39, 6).Value = 4532
$ws.Cells.Item(40, 6).Value = 5533
$ws.Cells.Item(43, 6).Value = 62
$ws.Cells.Item(44, 6).Value = 171
$ws.Cells.Item(45, 6).Value = 318
$ws.Cells.Item(47, 6).Value = 42
$ws.Cells.Item(48, 6).Value = 4110
$ws.Cells.Item(49, 6).Value = 125

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 4172
$ws.Cells.Item(4, 6).Value = 6
$ws.Cells.Item(13, 6).Value = 1029

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 450
$ws.Cells.Item(4, 6).Value = 75

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 450
$ws.Cells.Item(5, 6).Value = 75
$ws.Cells.Item(8, 6).Value = 901
$ws.Cells.Item(9, 6).Value = 1073
$ws.Cells.Item(10, 6).Value = 1548
$ws.Cells.Item(11, 6).Value = 330
$ws.Cells.Item(12, 6).Value = 671
$ws.Cells.Item(13, 6).Value = 12146
$ws.Cells.Item(14, 6).Value = 2164
$ws.Cells.Item(16, 6).Value = 1224
$ws.Cells.Item(17, 6).Value = 204
$ws.Cells.Item(18, 6).Value = 271
$ws.Cells.Item(19, 6).Value = 773
$ws.Cells.Item(20, 6).Value = 674
$ws.Cells.Item(22, 6).Value = 754
$ws.Cells.Item(23, 6).Value = 4065
$ws.Cells.Item(24, 6).Value = 4065
$ws.Cells.Item(25, 6).Value = 1114
$ws.Cells.Item(26, 6).Value = 860
$ws.Cells.Item(32, 6).Value = 1045
$ws.Cells.Item(33, 6).Value = 48
$ws.Cells.Item(34, 6).Value = 98
$ws.Cells.Item(37, 6).Value = 28
$ws.Cells.Item(38, 6).Value = 4415
$ws.Cells.Item(42, 6).Value = 171
$ws.Cells.Item(43, 6).Value = 318
$ws.Cells.Item(47, 6).Value = 4110

